$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend thick-bottom border row with two more empty, styled cells ---
$ws.Range("J3").Copy()
$ws.Range("K3:L3").PasteSpecial(-4122)

# --- Row 4: add year headers 2021 (K4) and 2022 (L4) ---
$ws.Range("J4").Copy()
$ws.Range("K4:L4").PasteSpecial(-4122)
$ws.Range("K4").Value = 2021
$ws.Range("L4").Value = 2022

# --- Row 5: "Всего" totals ---
$ws.Range("J5").Copy()
$ws.Range("K5:L5").PasteSpecial(-4122)
$ws.Range("K5:L5").HorizontalAlignment = 1
$ws.Range("K5").Value = 272.60000000000002
$ws.Range("L5").Value = 292.19961890663211

# --- Row 6: "в том числе:" (section header, empty values) ---
$ws.Range("J6").Copy()
$ws.Range("K6:L6").PasteSpecial(-4122)
$ws.Range("K6:L6").HorizontalAlignment = 1

# --- Row 7: "твердых" ---
$ws.Range("J7").Copy()
$ws.Range("K7:L7").PasteSpecial(-4122)
$ws.Range("K7:L7").HorizontalAlignment = 1
$ws.Range("K7").Value = 98.1
$ws.Range("L7").Value = 99.522498012012946

# --- Row 8: "газообразных и жидких" ---
$ws.Range("J8").Copy()
$ws.Range("K8:L8").PasteSpecial(-4122)
$ws.Range("K8:L8").HorizontalAlignment = 1
$ws.Range("K8").Value = 174.5
$ws.Range("L8").Value = 192.67712089461918

# --- Row 9: "из них:" (section header, empty values) ---
$ws.Range("J9").Copy()
$ws.Range("K9:L9").PasteSpecial(-4122)
$ws.Range("K9:L9").HorizontalAlignment = 1

# --- Row 10: "сернистого ангидрида" ---
$ws.Range("J10").Copy()
$ws.Range("K10:L10").PasteSpecial(-4122)
$ws.Range("K10:L10").HorizontalAlignment = 1
$ws.Range("K10").Value = 75.599999999999994
$ws.Range("L10").Value = 88.011952928467494

# --- Row 11: "окиси углерода" ---
$ws.Range("J11").Copy()
$ws.Range("K11:L11").PasteSpecial(-4122)
$ws.Range("K11:L11").HorizontalAlignment = 1
$ws.Range("K11").Value = 55.5
$ws.Range("L11").Value = 56.919430260413804

# --- Row 12: "окислов азота" (bottom thick-border row) ---
$ws.Range("J12").Copy()
$ws.Range("K12:L12").PasteSpecial(-4122)
$ws.Range("K12:L12").HorizontalAlignment = 1
$ws.Range("K12").Value = 24.9
$ws.Range("L12").Value = 24.176373211436804

# --- Update the selected cell to reflect final cursor position ---
$null = $ws.Range("N5").Select()
